# Scen_B_SYS_MaxGrowthRates.xlsx - RSD sheet update:
#  - Add a "maximum degrowth" row for Coal (after the existing Coal row)
#  - Add a "maximum degrowth" row for Peat (after the existing Peat row)
#  - Adjust a few growth-rate inputs (Biodiesel, Gas, Solar) in the data block
#  - Refresh the dependent TEXTJOIN summary formula

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSD")

# ---------------------------------------------------------------------
# 1. Insert two new rows into the UC table: one right after "Coal"
#    (row 9) for a Coal degrowth limit, and one right after "Peat"
#    (which will then sit at row 11) for a Peat degrowth limit.
#    Inserting rows here automatically re-points every formula below
#    (B/J/L/N columns, and the whole data block beneath) to the right
#    shifted cells, and extends the F6:F18 summary range to F6:F20.
# ---------------------------------------------------------------------
$ws.Rows("10:10").Insert()
$ws.Rows("12:12").Insert()

# ---------------------------------------------------------------------
# 2. Populate the new row 10 - "RSD maximum degrowth rate of Coal"
# ---------------------------------------------------------------------
$ws.Range("B10").Formula = '=_xlfn.TEXTJOIN("_",TRUE,"UC",A27,"MaxGrowth",B27)'
$ws.Range("C10").Formula = '=C9'
$ws.Range("F10").Value2 = "RSDCOA"
$ws.Range("G10").Formula = '=G9'
$ws.Range("H10").Value2 = 2021
$ws.Range("I10").Value2 = "UP"
$ws.Range("J10").Formula = '=1-C28'
$ws.Range("K10").Value2 = 1
$ws.Range("L10").Formula = '=-D28'
$ws.Range("M10").Value2 = 5
$ws.Range("N10").Formula = '=_xlfn.TEXTJOIN(" ",TRUE,A28, "maximum degrowth rate of",B27)'

# ---------------------------------------------------------------------
# 3. Populate the new row 12 - "RSD maximum degrowth rate of Peat"
# ---------------------------------------------------------------------
$ws.Range("B12").Formula = '=_xlfn.TEXTJOIN("_",TRUE,"UC",A29,"MaxGrowth",B28)'
$ws.Range("C12").Formula = '=C11'
$ws.Range("F12").Value2 = "RSDPEA"
$ws.Range("G12").Formula = '=G11'
$ws.Range("H12").Value2 = 2021
$ws.Range("I12").Value2 = "UP"
$ws.Range("J12").Formula = '=1-C28'
$ws.Range("K12").Value2 = 1
$ws.Range("L12").Formula = '=-D29'
$ws.Range("M12").Value2 = 5
$ws.Range("N12").Formula = '=_xlfn.TEXTJOIN(" ",TRUE,A29, "maximum degrowth rate of",B28)'

# ---------------------------------------------------------------------
# 4. Tweak the growth-rate inputs down in the data block (rows shifted
#    down by 2 after the inserts above: Biodiesel -> 25, Gas -> 31,
#    Solar -> 35).
# ---------------------------------------------------------------------
$ws.Range("C25").Value2 = 0.05   # Biodiesel max growth rate: 0.1 -> 0.05
$ws.Range("C31").Value2 = 0.01   # Gas max growth rate: 0.02 -> 0.01
$ws.Range("C35").Value2 = 0.1    # Solar max growth rate: 0.15 -> 0.1

# ---------------------------------------------------------------------
# 5. Update the view so the active selection matches the edited area.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("F26").Select()

$wb.Windows.Item(1).WindowState = -4143
